# Updated cryptos list - applies price/volume changes, a new "Frax" row
# insertion at row 34 (shifting rows 34-50 down by one), and replaces
# the former last two rows (Decentraland, Aave) with a single new
# "NEARProtocol" row at row 51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.622.07"
$ws.Range("E2").Value = "  +4.14%  "

$ws.Range("D3").Value = "1.747.74"
$ws.Range("E3").Value = "  +4.67%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.51"
$ws.Range("E5").Value = "  +3.62%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4815"
$ws.Range("E7").Value = "  +0.65%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2702"
$ws.Range("E8").Value = "  +2.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06260"
$ws.Range("E9").Value = "  +1.44%  "

$ws.Range("D10").Value = "1.746.70"
$ws.Range("E10").Value = "  +4.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07111"
$ws.Range("E11").Value = "  +1.67%  "

$ws.Range("E12").Value = "  +6.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6178"
$ws.Range("E13").Value = "  +4.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.513"
$ws.Range("E14").Value = "  +3.07%  "

$ws.Range("E15").Value = "  +2.71%  "

$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("D17").Value = "26.621.78"
$ws.Range("E17").Value = "  +4.19%  "

$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("E19").Value = "  +2.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.72"
$ws.Range("E20").Value = "  +2.44%  "

$ws.Range("D21").Value = "1.970.14"
$ws.Range("E21").Value = "  +4.55%  "

$ws.Range("E22").Value = "  +5.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.865"
$ws.Range("E23").Value = "  +1.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.356"
$ws.Range("E24").Value = "  +1.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.19"
$ws.Range("E25").Value = "  -0.32%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.48"
$ws.Range("E26").Value = "  +2.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.825"
$ws.Range("E27").Value = "  +6.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.414"
$ws.Range("E28").Value = "  +1.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.81"
$ws.Range("E29").Value = "  +2.94%  "

$ws.Range("E30").Value = "  +1.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.772"
$ws.Range("E31").Value = "  +3.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07900"
$ws.Range("E32").Value = "  +0.93%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04582"
$ws.Range("E33").Value = "  +8.54%  "

$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9998"
$ws.Range("E34").Value = "  +0.07%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.615"
$ws.Range("E35").Value = "  -0.25%  "

$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9992"
$ws.Range("E36").Value = "  +4.63%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6348"
$ws.Range("E37").Value = "  +4.18%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9591"
$ws.Range("E38").Value = "  +11.47%  "

$ws.Range("B39").Value = "Quant"
$ws.Range("C39").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "114.69"
$ws.Range("E39").Value = "  +18.79%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.472"
$ws.Range("E40").Value = "  -4.73%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.980"
$ws.Range("E41").Value = "  +5.86%  "

$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("E42").Value = "  +0.43%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01519"
$ws.Range("E43").Value = "  +2.90%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.718"
$ws.Range("E44").Value = "  +17.32%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3923"
$ws.Range("E45").Value = "  +4.16%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.753"
$ws.Range("E46").Value = "  +8.56%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1206"
$ws.Range("E47").Value = "  +8.03%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05329"
$ws.Range("E48").Value = "  +1.25%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.956"
$ws.Range("E49").Value = "  +7.89%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.88"
$ws.Range("E50").Value = "  +3.24%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.249"
$ws.Range("E51").Value = "  +4.02%  "
